$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 'running knee pads'
$ws.Range("A2").Value = 'youth yoga pants'
$ws.Range("A3").Value = 'boys lacrosse pants'
$ws.Range("A4").Value = 'knee length baseball pants'
$ws.Range("A5").Value = 'youth compression knee pad sleeve'
$ws.Range("A6").Value = 'yoga pad for knees'
$ws.Range("A7").Value = 'recovery leggings men'
$ws.Range("A8").Value = 'hex leg sleeves basketball'
$ws.Range("A9").Value = 'adult basketball'
$ws.Range("A10").Value = 'knee pads boys youth'
$ws.Range("A11").Value = 'wrestling leggings for men'
$ws.Range("A12").Value = 'softball sliding shorts youth'
$ws.Range("A13").Value = 'leggings sports'
$ws.Range("A14").Value = 'basketball padded shorts'
$ws.Range("A15").Value = 'best yoga knee pad'
$ws.Range("A16").Value = 'mens yoga pants tight'
$ws.Range("A17").Value = 'baseball tights'
$ws.Range("A18").Value = 'volleyball catcher'
$ws.Range("A19").Value = 'goalkeeper pads'
$ws.Range("A20").Value = 'youth baseball pants knee'
$ws.Range("A21").Value = 'paintball knee pads'
$ws.Range("A22").Value = 'youth compression tights boys'
$ws.Range("A23").Value = 'compression tights for boys basketball'
$ws.Range("A24").Value = 'breathable knee pads'
$ws.Range("A25").Value = 'mens spandex pants'
$ws.Range("A26").Value = 'soccer pads boys'
$ws.Range("A27").Value = 'hockey knee pads youth'
$ws.Range("A28").Value = 'padded shorts basketball'
$ws.Range("A29").Value = 'sliding shorts baseball youth'
$ws.Range("A30").Value = 'running pads'
$ws.Range("A31").Value = 'rodillera de basketball'
$ws.Range("A32").Value = 'leg tights for men'
$ws.Range("A33").Value = 'mens baseball pants black'
$ws.Range("A34").Value = 'sliding shorts boys'
$ws.Range("A35").Value = 'baseball pants knee'
$ws.Range("A36").Value = 'adult softball pants'
$ws.Range("A37").Value = 'catcher knee support'
$ws.Range("A38").Value = 'cycling pads for men'
$ws.Range("A39").Value = 'youth football girdle with knee pads'
$ws.Range("A40").Value = 'best knee pads for basketball'
$ws.Range("A41").Value = 'knee pads yoga'
$ws.Range("A42").Value = 'compression pads'
$ws.Range("A43").Value = 'knee pants boys'
$ws.Range("A44").Value = 'basketball chart'
$ws.Range("A45").Value = 'girl basketball knee pads'
$ws.Range("A46").Value = 'knee pads for basketball girls'
$ws.Range("A47").Value = 'knee pads men'
$ws.Range("A48").Value = 'basketball compression pants boys'
$ws.Range("A49").Value = 'mens compression knee'
$ws.Range("A50").Value = 'knee pads by design'
$ws.Range("A51").Value = 'good knee pads'
$ws.Range("A52").Value = 'wrestling knee pad youth'
$ws.Range("A53").Value = 'mens work pants with knee pad'
$ws.Range("A54").Value = 'kneepad pants'
$ws.Range("A55").Value = 'softball girls sliding shorts'
$ws.Range("A56").Value = 'men yoga capri'
$ws.Range("A57").Value = 'compression shorts lacrosse'
$ws.Range("A58").Value = 'black knee pads for basketball'
$ws.Range("A59").Value = 'mens 3/4 compression pants'
$ws.Range("A60").Value = 'softball sliding shorts youth girls'
$ws.Range("A61").Value = 'basketball knee pads girls'
$ws.Range("A62").Value = 'sliding baseball shorts'
$ws.Range("A63").Value = 'knee pads for big men'
$ws.Range("A64").Value = 'sleeve knee pads basketball'
$ws.Range("A65").Value = 'youth padded compression shorts basketball'
$ws.Range("A66").Value = 'mens compression recovery pants'
$ws.Range("A67").Value = 'arthritis friendly yoga'
$ws.Range("A68").Value = 'youth compression pants for boys'
$ws.Range("A69").Value = 'sports leggings for men'
$ws.Range("A70").Value = 'long knee pads basketball'
$ws.Range("A71").Value = 'best knee pads for yoga'
$ws.Range("A72").Value = 'rodillera basketball'
$ws.Range("A73").Value = 'basketball compression shorts'
$ws.Range("A74").Value = 'mountain biking knee pads for men'
$ws.Range("A75").Value = 'boys athletic tights youth'
$ws.Range("A76").Value = 'youth boys compression tights'
$ws.Range("A77").Value = 'boys compression pants football'
$ws.Range("A78").Value = 'padded sliding shorts mens'
$ws.Range("A79").Value = 'softball compression shorts'
$ws.Range("A80").Value = 'youth boys compression leggings'
$ws.Range("A81").Value = 'fitness gear floor guard'
$ws.Range("A82").Value = 'boys compression tights basketball'
$ws.Range("A83").Value = 'padded football pants youth'
$ws.Range("A84").Value = 'football girdle with knee pads'
$ws.Range("A85").Value = 'youth girls sliding shorts softball'
$ws.Range("A86").Value = 'construction need pads'
$ws.Range("A87").Value = 'softball sliding shorts for girls'
$ws.Range("A88").Value = 'softball pants men'
$ws.Range("A89").Value = 'pads for running'
$ws.Range("A90").Value = 'big and tall compression pants'
$ws.Range("A91").Value = 'hex pad knee sleeve'
$ws.Range("A92").Value = 'adult softball'
$ws.Range("A93").Value = 'compression knee sleeve men basketball'
$ws.Range("A94").Value = 'volleyball hand protector'
$ws.Range("A95").Value = 'yoga knee pad'
$ws.Range("A96").Value = 'knee sleeve wrestling youth'
$ws.Range("A97").Value = 'compression recovery pants'
$ws.Range("A98").Value = 'below knee leggings'
$ws.Range("A99").Value = 'basketball padded compression'
$ws.Range("A100").Value = 'knee length basketball shorts for women'
